$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999147246577,
    0.99813929396037304,
    0.9963348371025631,
    1.0008926035225949,
    0.99236616626994878,
    0.97196812360057039,
    0.97268710559477478,
    0.96914340740997562,
    0.96668249340082113,
    0.96516496445557876,
    0.96508571878414351,
    0.96519132663599361,
    0.94753227584138533,
    0.94057395315195969,
    0.93536037663740168,
    0.92974834953210883,
    0.91947536171702837,
    0.9152947458756191,
    0.99506382110156411,
    0.97011268813115081,
    0.96316716833555405,
    0.95161723522410024,
    0.98930091124916286,
    0.97628091515241011,
    0.96982407319404895,
    0.96831598768369709,
    0.96347365425195297,
    0.94201468175992531,
    0.92674851942911696,
    0.92017923963447923,
    0.91252640847662225,
    0.91084715798933291,
    0.91032716363825861
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
